# Apply the 2024-10-26 FlashScore refresh:
#  - reorder the "Odd_CS_3-3_HT" column so it sits right before "Odd_CS_0-1_HT"
#  - insert the new Adelaide United vs Central Coast Mariners fixture at the top
#  - refresh the odds that moved since the last scrape for the other fixtures
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Move column "Odd_CS_3-3_HT" (currently BC) so it lands right before "Odd_CS_0-1_HT" (AW) ---
$ws.Columns("BC").Cut() | Out-Null
$ws.Columns("AW").Insert() | Out-Null

# --- Insert the new match row at row 2, pushing the existing matches down one row ---
$ws.Rows(2).Insert() | Out-Null
# The insert copies the header row's bold/centered formatting onto the new row;
# strip it back to the plain (unstyled) look used by the other data rows.
$ws.Rows(2).ClearFormats() | Out-Null

# --- Fill in the new fixture: Adelaide United vs Central Coast Mariners ---
$ws.Range("A2").Value = "Yy9pYSf5"
$ws.Range("B2").Value = "26/10/2024"
$ws.Range("C2").Value = "03:00"
$ws.Range("D2").Value = "AUSTRALIA - A-LEAGUE"
$ws.Range("E2").Value = "Adelaide United"
$ws.Range("F2").Value = "Central Coast Mariners"
$ws.Range("G2").Value = 2.2
$ws.Range("H2").Value = 3.6
$ws.Range("I2").Value = 3.1
$ws.Range("J2").Value = 2.75
$ws.Range("K2").Value = 2.4
$ws.Range("L2").Value = 3.4
$ws.Range("M2").Value = 1.03
$ws.Range("N2").Value = 15
$ws.Range("O2").Value = 1.17
$ws.Range("P2").Value = 5
$ws.Range("Q2").Value = 1.57
$ws.Range("R2").Value = 2.38
$ws.Range("S2").Value = 1.29
$ws.Range("T2").Value = 3.5
$ws.Range("U2").Value = 1.5
$ws.Range("V2").Value = 2.5
$ws.Range("W2").Value = 12
$ws.Range("X2").Value = 13
$ws.Range("Y2").Value = 9
$ws.Range("Z2").Value = 21
$ws.Range("AA2").Value = 15
$ws.Range("AB2").Value = 21
$ws.Range("AC2").Value = 17
$ws.Range("AD2").Value = 7
$ws.Range("AE2").Value = 11
$ws.Range("AF2").Value = 34
$ws.Range("AG2").Value = 101
$ws.Range("AH2").Value = 15
$ws.Range("AI2").Value = 19
$ws.Range("AJ2").Value = 12
$ws.Range("AK2").Value = 34
$ws.Range("AL2").Value = 21
$ws.Range("AM2").Value = 23
$ws.Range("AN2").Value = 4.5
$ws.Range("AO2").Value = 11
$ws.Range("AP2").Value = 17
$ws.Range("AQ2").Value = 34
$ws.Range("AR2").Value = 41
$ws.Range("AS2").Value = 101
$ws.Range("AT2").Value = 3.5
$ws.Range("AU2").Value = 7
$ws.Range("AV2").Value = 41
$ws.Range("AW2").Value = 301
$ws.Range("AX2").Value = 5.5
$ws.Range("AY2").Value = 15
$ws.Range("AZ2").Value = 21
$ws.Range("BA2").Value = 51
$ws.Range("BB2").Value = 51
$ws.Range("BC2").Value = 101
$ws.Range("BD2").Value = 151

# --- Refresh odds for the other fixtures (values re-scraped since the last export) ---
# Row 3
$ws.Range("M3").Value = 1.05
$ws.Range("O3").Value = 1.25
# Row 4
$ws.Range("G4").Value = 4.33
$ws.Range("I4").Value = 1.7
$ws.Range("M4").Value = 1.02
$ws.Range("O4").Value = 1.15
$ws.Range("Q4").Value = 1.6
$ws.Range("R4").Value = 2.3
$ws.Range("S4").Value = 1.3
$ws.Range("T4").Value = 3.4
$ws.Range("W4").Value = 15
$ws.Range("X4").Value = 23
$ws.Range("Z4").Value = 41
$ws.Range("AH4").Value = 9
$ws.Range("AK4").Value = 15
$ws.Range("AL4").Value = 13
$ws.Range("AO4").Value = 21
$ws.Range("AT4").Value = 3.4
# Row 6
$ws.Range("H6").Value = 3.2
$ws.Range("K6").Value = 2.1
$ws.Range("L6").Value = 3.6
$ws.Range("M6").Value = 1.06
$ws.Range("N6").Value = 10
$ws.Range("O6").Value = 1.29
$ws.Range("P6").Value = 3.5
$ws.Range("Q6").Value = 2
$ws.Range("R6").Value = 1.85
$ws.Range("U6").Value = 1.69
$ws.Range("V6").Value = 2
$ws.Range("W6").Value = 8
$ws.Range("X6").Value = 11
$ws.Range("AC6").Value = 10
$ws.Range("AD6").Value = 6
$ws.Range("AJ6").Value = 12
$ws.Range("AL6").Value = 26
$ws.Range("AM6").Value = 34
$ws.Range("AU6").Value = 8
$ws.Range("AZ6").Value = 26
$ws.Range("BB6").Value = 81
